# Apply F-column ("想去人数" / interest-count) updates per the commit
# "Update gh-pages to output generated at 456a3b4" across all 4 sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1906
$ws.Range("F4").Value = 80
$ws.Range("F5").Value = 806
$ws.Range("F6").Value = 129
$ws.Range("F7").Value = 523
$ws.Range("F8").Value = 991
$ws.Range("F10").Value = 1307
$ws.Range("F11").Value = 1599
$ws.Range("F12").Value = 79
$ws.Range("F13").Value = 1603
$ws.Range("F14").Value = 361
$ws.Range("F15").Value = 1730
$ws.Range("F16").Value = 821
$ws.Range("F17").Value = 1170
$ws.Range("F19").Value = 2044
$ws.Range("F21").Value = 837
$ws.Range("F22").Value = 1028
$ws.Range("F23").Value = 563
$ws.Range("F25").Value = 1333
$ws.Range("F26").Value = 1116
$ws.Range("F27").Value = 105
$ws.Range("F28").Value = 583
$ws.Range("F29").Value = 1236
$ws.Range("F31").Value = 1222
$ws.Range("F32").Value = 63
$ws.Range("F33").Value = 1170
$ws.Range("F34").Value = 334
$ws.Range("F35").Value = 93
$ws.Range("F38").Value = 1742
$ws.Range("F39").Value = 396
$ws.Range("F40").Value = 19
$ws.Range("F44").Value = 852
$ws.Range("F45").Value = 820
$ws.Range("F47").Value = 816
$ws.Range("F48").Value = 127

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value = 2636
$ws.Range("F13").Value = 45
$ws.Range("F14").Value = 82
$ws.Range("F21").Value = 95198
$ws.Range("F25").Value = 43
$ws.Range("F26").Value = 43
$ws.Range("F27").Value = 199
$ws.Range("F28").Value = 262
$ws.Range("F30").Value = 240
$ws.Range("F32").Value = 62
$ws.Range("F37").Value = 195
$ws.Range("F41").Value = 74
$ws.Range("F43").Value = 151

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 277
$ws.Range("F5").Value = 2960
$ws.Range("F6").Value = 4738
$ws.Range("F9").Value = 616
$ws.Range("F10").Value = 824
$ws.Range("F11").Value = 493
$ws.Range("F12").Value = 480
$ws.Range("F13").Value = 1203
$ws.Range("F14").Value = 335
$ws.Range("F15").Value = 855

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1906
$ws.Range("F3").Value = 277
$ws.Range("F5").Value = 4738
$ws.Range("F6").Value = 616
$ws.Range("F7").Value = 824
$ws.Range("F8").Value = 493
$ws.Range("F9").Value = 480
$ws.Range("F10").Value = 480
$ws.Range("F11").Value = 1203
$ws.Range("F12").Value = 523
$ws.Range("F13").Value = 991
$ws.Range("F15").Value = 1307
$ws.Range("F16").Value = 1599
$ws.Range("F17").Value = 79
$ws.Range("F18").Value = 1603
$ws.Range("F20").Value = 82
$ws.Range("F21").Value = 1730
$ws.Range("F22").Value = 1170
$ws.Range("F23").Value = 855
$ws.Range("F24").Value = 855
$ws.Range("F25").Value = 2044
$ws.Range("F27").Value = 837
$ws.Range("F28").Value = 1028
$ws.Range("F29").Value = 563
$ws.Range("F30").Value = 1333
$ws.Range("F32").Value = 1116
$ws.Range("F33").Value = 105
$ws.Range("F34").Value = 1236
$ws.Range("F36").Value = 1222
$ws.Range("F37").Value = 63
$ws.Range("F39").Value = 43
$ws.Range("F40").Value = 1170
$ws.Range("F41").Value = 334
$ws.Range("F44").Value = 1742
$ws.Range("F45").Value = 19
$ws.Range("F48").Value = 852
$ws.Range("F49").Value = 821
$ws.Range("F50").Value = 816
$ws.Range("F51").Value = 127
